$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.926.78'
$ws.Range("E2").Value = '  -3.07%  '

$ws.Range("D3").Value = '3.310.54'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''573.44'
$ws.Range("E5").Value = '  -1.72%  '

$ws.Range("D6").Value = '''180.91'
$ws.Range("E6").Value = '  -3.42%  '

$ws.Range("D7").Value = '''0.616'
$ws.Range("E7").Value = '  +2.92%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  -2.15%  '

$ws.Range("D10").Value = '''6.65'
$ws.Range("E10").Value = '  -0.53%  '

$ws.Range("E11").Value = '  -1.74%  '

$ws.Range("D12").Value = '3.887.26'
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("E13").Value = '  -1.34%  '

$ws.Range("D14").Value = '''26.66'
$ws.Range("E14").Value = '  -3.69%  '

$ws.Range("D15").Value = '66.063.35'
$ws.Range("E15").Value = '  -3.07%  '

$ws.Range("E16").Value = '  -1.21%  '

$ws.Range("D17").Value = '3.283.43'
$ws.Range("E17").Value = '  -1.03%  '

$ws.Range("D18").Value = '''436.41'
$ws.Range("E18").Value = '  -2.30%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''13.52'
$ws.Range("E19").Value = '  -0.59%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''5.65'
$ws.Range("E20").Value = '  -1.65%  '

$ws.Range("D21").Value = '''7.55'
$ws.Range("E21").Value = '  -2.31%  '

$ws.Range("D22").Value = '''73.16'
$ws.Range("E22").Value = '  -3.03%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.54%  '

$ws.Range("D24").Value = '''0.519'
$ws.Range("E24").Value = '  +0.31%  '

$ws.Range("E25").Value = '  -2.88%  '

$ws.Range("D26").Value = '''0.193'
$ws.Range("E26").Value = '  +2.17%  '

$ws.Range("D27").Value = '''9.07'
$ws.Range("E27").Value = '  -1.53%  '

$ws.Range("E28").Value = '  -0.35%  '

$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("D30").Value = '''22.66'
$ws.Range("E30").Value = '  -1.48%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").Value = '''5.23'
$ws.Range("E32").Value = '  -2.90%  '

$ws.Range("D33").Value = '''6.75'
$ws.Range("E33").Value = '  -1.06%  '

$ws.Range("E34").Value = '  -3.12%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''1.48'
$ws.Range("E35").Value = '  -3.34%  '

$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").Value = '''159.76'
$ws.Range("E36").Value = '  -2.47%  '

$ws.Range("D37").Value = '''27.68'
$ws.Range("E37").Value = '  +2.35%  '

$ws.Range("E38").Value = '  -5.38%  '

$ws.Range("D39").Value = '2.829.40'
$ws.Range("E39").Value = '  +4.83%  '

$ws.Range("D40").Value = '''0.788'
$ws.Range("E40").Value = '  -0.36%  '

$ws.Range("E41").Value = '  -2.67%  '

$ws.Range("D42").Value = '''6.18'
$ws.Range("E42").Value = '  -3.74%  '

$ws.Range("D43").Value = '''40.30'
$ws.Range("E43").Value = '  -1.23%  '

$ws.Range("D44").Value = '''0.0665'
$ws.Range("E44").Value = '  -1.85%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.35'
$ws.Range("E45").Value = '  -3.10%  '

$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''24.11'
$ws.Range("E46").Value = '  -2.33%  '

$ws.Range("D47").Value = '''325.17'
$ws.Range("E47").Value = '  -0.81%  '

$ws.Range("E48").Value = '  -2.18%  '

$ws.Range("E49").Value = '  +1.49%  '

$ws.Range("E50").Value = '  -1.55%  '

$ws.Range("E51").Value = '  -2.09%  '

